$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new reference entry in C12: "Channel: 3Blue1Brown"
$ws.Range("C12").Value = "Channel: 3Blue1Brown"
$ws.Range("C12").WrapText = $true

# Update selection / view: select G14
$ws.Range("G14").Select()
